# Re-running the PanelApp query later the same day: the "data" sheet's
# per-row query timestamps (column F) are refreshed, and a new "metadata"
# tab summarising the query itself is appended after "data".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Refresh the per-gene query timestamps on the "data" sheet ---
$ws.Range("F2").Value  = "2021-10-05 14:21:08.332421"
$ws.Range("F3").Value  = "2021-10-05 14:21:08.332429"
$ws.Range("F4").Value  = "2021-10-05 14:21:08.332433"
$ws.Range("F5").Value  = "2021-10-05 14:21:08.332435"
$ws.Range("F6").Value  = "2021-10-05 14:21:08.332438"
$ws.Range("F7").Value  = "2021-10-05 14:21:08.332441"
$ws.Range("F8").Value  = "2021-10-05 14:21:08.332444"
$ws.Range("F9").Value  = "2021-10-05 14:21:08.332446"
$ws.Range("F10").Value = "2021-10-05 14:21:08.332449"
$ws.Range("F11").Value = "2021-10-05 14:21:08.332452"
$ws.Range("F12").Value = "2021-10-05 14:21:08.332455"
$ws.Range("F13").Value = "2021-10-05 14:21:08.332457"
$ws.Range("F14").Value = "2021-10-05 14:21:08.332459"
$ws.Range("F15").Value = "2021-10-05 14:21:08.332462"
$ws.Range("F16").Value = "2021-10-05 14:21:08.332465"
$ws.Range("F17").Value = "2021-10-05 14:21:08.332467"
$ws.Range("F18").Value = "2021-10-05 14:21:08.332470"
$ws.Range("F19").Value = "2021-10-05 14:21:08.332473"
$ws.Range("F20").Value = "2021-10-05 14:21:08.332475"
$ws.Range("F21").Value = "2021-10-05 14:21:08.332478"
$ws.Range("F22").Value = "2021-10-05 14:21:08.332480"
$ws.Range("F23").Value = "2021-10-05 14:21:08.332483"
$ws.Range("F24").Value = "2021-10-05 14:21:08.332485"
$ws.Range("F25").Value = "2021-10-05 14:21:08.332488"
$ws.Range("F26").Value = "2021-10-05 14:21:08.332491"
$ws.Range("F27").Value = "2021-10-05 14:21:08.332494"

# --- Add the "metadata" tab right after "data" ---
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$new.Name = "metadata"

# Header row
$new.Range("B1").Value = "data_name"
$new.Range("C1").Value = "data_id"
$new.Range("D1").Value = "data_version"
$new.Range("E1").Value = "data_version_created"
$new.Range("F1").Value = "panel_query_time"
$new.Range("G1").Value = "panel_get_request"

# Match the bold/bordered header styling already used on the "data" sheet
$ws.Range("B1:F1").Copy()
$new.Range("B1:F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B1").Copy()
$new.Range("G1").PasteSpecial(-4122)      # xlPasteFormats

# Data row
$new.Range("A2").Value = 0
$new.Range("B2").Value = "Inherited ovarian cancer (without breast cancer)"
$new.Range("C2").Value = 143

# Keep "2.21" as text (not the number 2.21)
$new.Range("D2").NumberFormat = "@"
$new.Range("D2").Value = "2.21"
$new.Range("D2").Style = "Normal"

$new.Range("E2").Value = "2021-07-15T09:17:52.140378Z"
$new.Range("F2").Value = "2021-10-05 14:21:08.329288"
$new.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/143/?format=json"

# A2's index style also matches the "data" sheet's index-column styling
$ws.Range("A2").Copy()
$new.Range("A2").PasteSpecial(-4122)      # xlPasteFormats
